$d = $word.ActiveDocument

# Locate the "KHOA CNTT&TT" run on the cover page and split it into two
# runs: "KHOA " and the spelled-out "CÔNG NGHỆ THÔNG TIN VÀ TRUYỀN THÔNG",
# both keeping the original bold/size formatting.
$finder = $d.Content
$found = $finder.Find.Execute("KHOA CNTT&TT", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $target = $d.Range($finder.Start, $finder.End)

    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r>
<w:rPr><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>
<w:t xml:space="preserve">KHOA </w:t>
</w:r>
<w:r>
<w:rPr><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>
<w:t>CÔNG NGHỆ THÔNG TIN VÀ TRUYỀN THÔNG</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

    $target.InsertXML($xml)
}
